$data = @(
  @("Crafty", 354, "Active", 17, 4, 11, 0, 11, 56, 3.29, 14.01, "03. 4 of a Kind", "2,", 14.03),
  @("Illya", 355, "Active", 27, 5, 10, 0, 10, 109, 4.04, 4.09, $null, $null, $null),
  @("Jon", 357, "Active", 164, 17, 13, 0, 13, 573, 3.49, 13.08, "03. 4 of a Kind", "Ace,", 12.04),
  @("Maisy", 360, "Active", 96, 11, 10, 0, 10, 335, 3.49, 16.03, "02. Straight Flush", "6,7,8,9,10,(D)", 15.09),
  @("Mark", 361, "Active", 117, 13, 9, 0, 9, 428, 3.66, 14.12, "03. 4 of a Kind", "Ace,", 14.02),
  @("Matt", 362, "Active", 159, 17, 11, 0, 11, 589, 3.7, 13.12, "03. 4 of a Kind", "King,", 16.02),
  @("Pepe", 364, "Active", 77, 10, 9, 1, 10, 269, 3.49, 17.02, "03. 4 of a Kind", "Ace,", 15.08),
  @("Prashant", 365, "Active", 12, 3, 9, 0, 9, 46, 3.83, 16.01, "03. 4 of a Kind", "King,", 16.01),
  @("Richard", 366, "Active", 114, 17, 12, 0, 12, 489, 4.29, 12.05, "04. Full House (3,2)", "Ace,9,10,Jack,(9)", 15.03),
  @("Alex", 348, "Active", 72, 12, 9, 0, 9, 246, 3.42, 8.07, "03. 4 of a Kind", "King,", 11.02),
  @("Andy", 349, "Active", 169, 17, 11, 0, 11, 714, 4.22, 11.04, "03. 4 of a Kind", "King,", 16.05),
  @("Anthony", 350, "Active", 96, 11, 10, 0, 10, 380, 3.96, 9.08, "02. Straight Flush", "7,8,9,10,Jack,(H)", 12.12),
  @("Bob", 351, "Inactive", 35, 9, 11, 0, 11, 139, 3.97, 7.01, "04. Full House (3,2)", "4,5,", 10.12),
  @("Chris", 353, "Inactive", 1, 1, 3, 0, 3, 3, 3, 5.05, $null, $null, $null),
  @("John", 113, "Inactive", 9, 3, 12, 0, 12, 46, 5.11, 5.06, $null, $null, $null),
  @("Richard Snr", 367, "Inactive", 1, 1, 1, 0, 1, 1, 1, 3.11, $null, $null, $null),
  @("Stuart", 368, "Inactive", 48, 7, 12, 0, 12, 183, 3.81, 4.02, $null, $null, $null),
  @("Matthew", 363, "Inactive", 2, 1, 3, 0, 3, 5, 2.5, 4.06, $null, $null, $null),
  @("Jonathan", 358, "Inactive", 10, 3, 7, 0, 7, 35, 3.5, 4.07, $null, $null, $null),
  @("Keith", 359, "Inactive", 98, 12, 10, 0, 10, 333, 3.4, 13.04, "03. 4 of a Kind", "8,", 12.06),
  @("Jim", 356, "Inactive", 3, 3, 9, 0, 9, 17, 5.67, 7.05, $null, $null, $null),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    if ($null -eq $val) {
      $ws.Cells.Item($r, $c).ClearContents()
    } else {
      $ws.Cells.Item($r, $c).Value = $val
    }
    $c = $c + 1
  }
  $r = $r + 1
}

$ws.Range("Q7").Select()
